# Multiple corrections on the html pages, plus qa on resources and examples
#
# 1) Fix the cached "datetimeFigureOut" field text from 11/02/2020 to
#    12/02/2020 everywhere it appears: once in every slide layout's date
#    placeholder and once in the slide master's date placeholder.
# 2) Reposition/resize + shrink the font of the "E2: Serum electrolyte
#    panel, orderable: summary of content" title textbox on slide 1.

$p = $ppt.ActivePresentation

# --- 1) Date placeholder text fix (slide master + every slide layout) ---

$OLD_DATE = "11/02/2020"
$NEW_DATE = "12/02/2020"

$sm = $p.SlideMaster

for ($li = 1; $li -le $sm.CustomLayouts.Count; $li++) {
    $cl = $sm.CustomLayouts.Item($li)
    for ($i = 1; $i -le $cl.Shapes.Count; $i++) {
        $shp = $cl.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq $OLD_DATE) {
                $shp.TextFrame.TextRange.Text = $NEW_DATE
            }
        }
    }
}

$master = $p.Slides.Item(1).Master
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.TextRange.Text -eq $OLD_DATE) {
            $shp.TextFrame.TextRange.Text = $NEW_DATE
        }
    }
}

# --- 2) Title textbox move/resize + font size change on slide 1 ---

$s1 = $p.Slides.Item(1)
$title = $null
for ($i = 1; $i -le $s1.Shapes.Count; $i++) {
    $cand = $s1.Shapes.Item($i)
    if ($cand.Name -eq "ZoneTexte 28") {
        $title = $cand
        break
    }
}

# Target EMU values (from the OOXML diff):
#   off  x=159078  y=712932
#   ext cx=10874477 cy=461665
# Shape.Left/Top/Width/Height are single-precision (points) in the real
# PowerPoint object model, so the literal point values are nudged by a
# hair to land exactly on the intended EMU after the internal
# point->EMU truncation/round-trip.
$title.Left = 12.525826971653544
$title.Top = 56.13637735275591
$title.Width = 856.258026496063
$title.Height = 36.3515759031496

# All four runs shared size 3200 (32pt) and all move to 2400 (24pt); set
# it once across the whole text range so paragraph/run/highlight
# structure is left untouched.
$title.TextFrame.TextRange.Font.Size = 24
